$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.016.11'
$ws.Range("E2").Value = '  -1.55%  '
$ws.Range("D3").Value = '1.978.50'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.014'
$ws.Range("E4").Value = '  +0.57%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.37'
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.012'
$ws.Range("E6").Value = '  +0.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4971'
$ws.Range("E7").Value = '  -0.66%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4194'
$ws.Range("E8").Value = '  -0.50%  '
$ws.Range("E9").Value = '  +4.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09255'
$ws.Range("E10").Value = '  +4.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.094'
$ws.Range("E11").Value = '  -2.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.77'
$ws.Range("E12").Value = '  -2.29%  '
$ws.Range("D13").Value = '2.001.45'
$ws.Range("E13").Value = '  -1.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.921'
$ws.Range("E14").Value = '  -2.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.441'
$ws.Range("E15").Value = '  -1.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.014'
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001109'
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.10'
$ws.Range("E18").Value = '  -4.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06740'
$ws.Range("E19").Value = '  +1.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.23'
$ws.Range("E20").Value = '  -2.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.012'
$ws.Range("E21").Value = '  +0.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.968'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").Value = '29.064.49'
$ws.Range("E23").Value = '  -1.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.97'
$ws.Range("E24").Value = '  +0.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.287'
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").Value = '2.249.51'
$ws.Range("E26").Value = '  -0.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.80'
$ws.Range("E27").Value = '  +0.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '156.30'
$ws.Range("E28").Value = '  -1.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.285'
$ws.Range("E29").Value = '  -3.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.270'
$ws.Range("E30").Value = '  -2.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '127.54'
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.050'
$ws.Range("E32").Value = '  +0.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09828'
$ws.Range("E33").Value = '  -1.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.529'
$ws.Range("E34").Value = '  -1.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.824'
$ws.Range("E35").Value = '  -0.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.742'
$ws.Range("E36").Value = '  -1.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02423'
$ws.Range("E37").Value = '  -1.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.316'
$ws.Range("E38").Value = '  +2.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06398'
$ws.Range("E39").Value = '  +0.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.041'
$ws.Range("E40").Value = '  -5.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6484'
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.48'
$ws.Range("E42").Value = '  -2.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.2000'
$ws.Range("E43").Value = '  -3.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.012'
$ws.Range("E44").Value = '  +0.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6209'
$ws.Range("E45").Value = '  -2.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.357'
$ws.Range("E46").Value = '  +6.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '13.32'
$ws.Range("E47").Value = '  -0.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.182'
$ws.Range("E48").Value = '  -0.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.489'
$ws.Range("E49").Value = '  -0.72%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000330'
$ws.Range("E50").Value = '  -1.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06969'
$ws.Range("E51").Value = '  -0.60%  '
